$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new time entries on row 35 (new shift), continuing the pattern
# from prior rows: start/end time with shared formula duration in E.
$ws.Range("C35").Value = 0.9375
$ws.Range("D35").Value = 0.99930555555555556
$ws.Range("C35").NumberFormat = "h:mm"
$ws.Range("D35").NumberFormat = "h:mm:ss"
$ws.Range("E35").Formula = "=D35-C35"

# Move the active selection to F39, matching the updated view state.
$ws.Range("F39").Select()
